# Data consolidation on the "Stats" sheet: rename a couple of stats so their
# two-letter symbols no longer collide with another stat's symbol.
#
#   - "Attack Power" (AP) is renamed to "Impact" (IM), freeing up "AP" to be
#     the (now unique) symbol for "Appearance".
#   - "Vision" (VI) is renamed to "Sight" (SI), freeing up "VI" to be the
#     (now unique) symbol for "Vitality".
#   - A brand new stat "Structure" (ST) is inserted between "Opacity" and
#     "Life Force".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(8)   # "Stats" sheet

# --- Row 3: Attack Power / AP -> Impact / IM -----------------------------
$ws.Range("A3").Value2 = "Impact"
$ws.Range("B3").Value2 = "IM"

# --- Row 14: Appearance now owns its own "AP" symbol ----------------------
# (previously shown in maroon because it reused the "Attack Power" symbol;
# that conflict is gone now, so restore the normal black text color)
$ws.Range("B14").Font.Color = 0

# --- Row 17: Vision / VI -> Sight / SI ------------------------------------
$ws.Range("A17").Value2 = "Sight"
$ws.Range("B17").Value2 = "SI"

# --- Insert a new row 22 for "Structure" / "ST" ---------------------------
$ws.Rows.Item(22).Insert()
$ws.Range("A22").Value2 = "Structure"
$ws.Range("B22").Value2 = "ST"
$ws.Range("A22:B22").Font.Color = 0

# --- Row 24 (the old "Vitality" row, shifted down by the insert above):
# it now owns its own "VI" symbol, so clear the old maroon conflict marker.
$ws.Range("B24").Font.Color = 0

# --- Update the active selection to match the author's final cursor position
$ws.Activate()
$ws.Range("C17").Select()
